# Update countries & provincias Spain
# Applies the diff between before.xlsx and the target workbook:
#   1) Three adjacent shared-string reorderings (country names), which in the
#      underlying case-count-sorted table correspond to rows swapping places:
#        - Eslovaquia / Mozambique (rows 110/111)
#        - Santa Lucia / Timor Oriental (rows 204/205 - values already equal)
#        - Montserrat / Islas Malvinas (rows 214/215)
#   2) Refreshed COVID-19 statistics for a number of countries.
#   3) Updated "last refreshed" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Re-ordered country names (rows keep their position, the country that
#        occupies the row changes because the pair was re-sorted) ---

# Eslovaquia now sits above Mozambique (row 110 / row 111)
$ws.Range("A110").Value = "Eslovaquia"
$ws.Range("A111").Value = "Mozambique"

# Santa Lucia now sits above Timor Oriental (row 204 / row 205)
$ws.Range("A204").Value = "Santa Lucia"
$ws.Range("A205").Value = "Timor Oriental"

# Montserrat now sits above Islas Malvinas (row 214 / row 215)
$ws.Range("A214").Value = "Montserrat"
$ws.Range("A215").Value = "Islas Malvinas"

# --- 2) Updated statistics ---

# Row 5 (Peru)
$ws.Range("B5").Value = 5312537
$ws.Range("C5").Value = 7062
$ws.Range("E5").Value = 1018456
$ws.Range("G5").Value = 25
$ws.Range("H5").Value = 85650

# Row 16
$ws.Range("B16").Value = 419043
$ws.Range("C16").Value = 2845
$ws.Range("D16").Value = 357632
$ws.Range("E16").Value = 37293
$ws.Range("G16").Value = 166
$ws.Range("H16").Value = 24118

# Row 33
$ws.Range("B33").Value = 111550
$ws.Range("C33").Value = 1333
$ws.Range("D33").Value = 89119
$ws.Range("E33").Value = 18029
$ws.Range("G33").Value = 42
$ws.Range("H33").Value = 4402

# Row 59
$ws.Range("B59").Value = 50634
$ws.Range("C59").Value = 381
$ws.Range("D59").Value = 46645
$ws.Range("E59").Value = 3566
$ws.Range("G59").Value = 4
$ws.Range("H59").Value = 423

# Row 64
$ws.Range("B64").Value = 45857
$ws.Range("C64").Value = 97
$ws.Range("D64").Value = 45029
$ws.Range("E64").Value = 533

# Row 95
$ws.Range("B95").Value = 10488
$ws.Range("C95").Value = 32
$ws.Range("D95").Value = 9891
$ws.Range("E95").Value = 329

# Row 110 (now Eslovaquia)
$ws.Range("B110").Value = 6546
$ws.Range("C110").Value = 290
$ws.Range("D110").Value = 3519
$ws.Range("E110").Value = 2988
$ws.Range("H110").Value = 39

# Row 111 (now Mozambique)
$ws.Range("B111").Value = 6264
$ws.Range("D111").Value = 3502
$ws.Range("E111").Value = 2722
$ws.Range("H111").Value = 40

# Row 140
$ws.Range("D140").Value = 3070
$ws.Range("E140").Value = 198

# Row 145
$ws.Range("B145").Value = 2699
$ws.Range("C145").Value = 65
$ws.Range("D145").Value = 2017
$ws.Range("E145").Value = 663
$ws.Range("G145").Value = 2
$ws.Range("H145").Value = 19

# Row 182
$ws.Range("B182").Value = 350
$ws.Range("C182").Value = 4
$ws.Range("D182").Value = 322
$ws.Range("E182").Value = 28

# Row 214 (now Montserrat)
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

# Row 215 (now Islas Malvinas)
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0

# --- 3) Updated timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Septiembre de 2020 a las 13:10"
